$d = $word.ActiveDocument

# Locate the end of the document's last paragraph (the "Contact ..." paragraph)
# and collapse the range there so nothing existing gets disturbed.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)

# Body text for the new "Possible Destin Meetings" paragraph.
$heading = "Possible Destin Meetings"
$body = "  The Destin Library has nice facilities for small groups of  participants in a computer class.  Anyone interested in attending such a class for Linux please contact Tom."

# Insert: <paragraph break> + heading + body + <paragraph break>
# Leading CR ends the current (last) paragraph without touching its existing
# runs; the trailing CR creates the final, empty paragraph with no run.
$r.Text = [char]13 + $heading + $body + [char]13

# Grab the freshly created "Possible Destin Meetings" paragraph (second to
# last paragraph now; the very last is the new trailing empty paragraph).
$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)
$newRange = $newPara.Range

# Bold + underline just the heading portion of the new paragraph's text.
$headingRange = $d.Range($newRange.Start, $newRange.Start + $heading.Length)
$headingRange.Bold = 1
$headingRange.Font.Underline = 1
